# Update the ESPN comparison table with the latest NBA games (Saturday 10th Feb 2024),
# replacing the previous Friday 9th Feb 2024 slate.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "NBA, Saturday 10th Feb 2024"
$ws.Range("B1").Value = "Ballgorithm"
$ws.Range("C1").Value = "ESPN"

$ws.Range("A2").Value = "Oklahoma City Thunder (35-16) vs Dallas Mavericks (29-23)"
$ws.Range("B2").Value = "Oklahoma City Thunder (76.92%)"
$ws.Range("C2").Value = "Oklahoma City Thunder (64.5%)"

$ws.Range("A3").Value = "Detroit Pistons (8-43) vs Los Angeles Clippers (34-16)"
$ws.Range("B3").Value = "Los Angeles Clippers (79.17%)"
$ws.Range("C3").Value = "Los Angeles Clippers (94.0%)"

$ws.Range("A4").Value = "San Antonio Spurs (10-42) vs Brooklyn Nets (20-31)"
$ws.Range("B4").Value = "San Antonio Spurs (42.86%)"
$ws.Range("C4").Value = "Brooklyn Nets (74.8%)"

$ws.Range("A5").Value = "Memphis Grizzlies (18-34) vs Charlotte Hornets (10-41)"
$ws.Range("B5").Value = "Memphis Grizzlies (20.83%)"
$ws.Range("C5").Value = "Memphis Grizzlies (66.9%)"

$ws.Range("A6").Value = "Chicago Bulls (25-27) vs Orlando Magic (28-24)"
$ws.Range("B6").Value = "Orlando Magic (69.57%)"
$ws.Range("C6").Value = "Orlando Magic (68.5%)"

$ws.Range("A7").Value = "Philadelphia 76ers (30-21) vs Washington Wizards (9-42)"
$ws.Range("B7").Value = "Philadelphia 76ers (62.96%)"
$ws.Range("C7").Value = "Philadelphia 76ers (63.6%)"

$ws.Range("A8").Value = "Houston Rockets (23-28) vs Atlanta Hawks (23-29)"
$ws.Range("B8").Value = "Houston Rockets (66.67%)"
$ws.Range("C8").Value = "Houston Rockets (56.4%)"

$ws.Range("A9").Value = "Indiana Pacers (29-25) vs New York Knicks (33-19)"
$ws.Range("B9").Value = "New York Knicks (73.08%)"
$ws.Range("C9").Value = "New York Knicks (56.0%)"

$ws.Range("A10").Value = "Cleveland Cavaliers (34-16) vs Toronto Raptors (19-33)"
$ws.Range("B10").Value = "Cleveland Cavaliers (69.23%)"
$ws.Range("C10").Value = "Cleveland Cavaliers (78.3%)"

$ws.Range("A11").Value = "Phoenix Suns (31-21) vs Golden State Warriors (24-25)"
$ws.Range("B11").Value = "Phoenix Suns (59.26%)"
$ws.Range("C11").Value = "Golden State Warriors (55.2%)"

$ws.Range("A12").Value = "New Orleans Pelicans (30-22) vs Portland Trail Blazers (15-36)"
$ws.Range("B12").Value = "New Orleans Pelicans (60.00%)"
$ws.Range("C12").Value = "New Orleans Pelicans (81.8%)"

# Reset the active cell/selection back to the top-left of the sheet.
$ws.Range("A1").Select()
